# New crime data collected - update 68th Precinct weekly CompStat report
# (volume/week header text, and the weekly crime-complaint figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header strings -------------------------------------------------
# "Volume 30   Number  29"  ->  "Volume 30   Number  30"
$ws.Range("A8").Value = "Volume 30   Number  30"
# "Report Covering the Week  7/17/2023  Through  7/23/2023"
#   -> "Report Covering the Week  7/24/2023  Through  7/30/2023"
$ws.Range("C9").Value = "Report Covering the Week  7/24/2023  Through  7/30/2023"

# --- Row 15 -----------------------------------------------------------
$ws.Range("N15").Value = -62.5

# --- Row 16 -------------------------------------------------------------
$ws.Range("C16").Value = "'0"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null

$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 43
$ws.Range("J16").Value = 29
$ws.Range("K16").Value = 48.275862068965
$ws.Range("L16").Value = 48.275862068965
$ws.Range("M16").Value = -24.561403508771
$ws.Range("N16").Value = -87.125748502994

# --- Row 17 -------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -38.461538461538
$ws.Range("I17").Value = 86
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 43.333333333333
$ws.Range("L17").Value = 65.384615384615
$ws.Range("M17").Value = 45.762711864406
$ws.Range("N17").Value = -47.878787878787

# --- Row 18 -------------------------------------------------------------
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 60
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = 27.659574468085
$ws.Range("L18").Value = -11.764705882352
$ws.Range("M18").Value = -57.446808510638
$ws.Range("N18").Value = -90.4

# --- Row 19 -------------------------------------------------------------
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -18.181818181818
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = -12.195121951219
$ws.Range("I19").Value = 253
$ws.Range("J19").Value = 319
$ws.Range("K19").Value = -20.689655172413
$ws.Range("L19").Value = 21.634615384615
$ws.Range("M19").Value = 45.402298850574
$ws.Range("N19").Value = -7.326007326007

# --- Row 20 -------------------------------------------------------------
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -26.666666666666
$ws.Range("I20").Value = 68
$ws.Range("J20").Value = 72
$ws.Range("K20").Value = -5.555555555555
$ws.Range("L20").Value = 74.358974358974
$ws.Range("M20").Value = -27.659574468085
$ws.Range("N20").Value = -94.092093831450

# --- Row 21 (precinct total) --------------------------------------------
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 11.764705882352
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = -23.456790123456
$ws.Range("I21").Value = 518
$ws.Range("J21").Value = 536
$ws.Range("K21").Value = -3.358208955223
$ws.Range("L21").Value = 28.535980148883
$ws.Range("M21").Value = -2.631578947368
$ws.Range("N21").Value = -79.820802493182

# --- Row 22 (Transit) ----------------------------------------------------
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("C23:D23").Copy() | Out-Null
$ws.Range("C22:D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null

# --- Row 24 (Petit Larceny) ----------------------------------------------
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = -42.857142857142
$ws.Range("F24").Value = 120
$ws.Range("H24").Value = -27.272727272727
$ws.Range("I24").Value = 883
$ws.Range("J24").Value = 1076
$ws.Range("K24").Value = -17.936802973977
$ws.Range("L24").Value = 61.721611721611
$ws.Range("M24").Value = 31.398809523809

# --- Row 25 (Misd. Assault) -----------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -7.692307692307
$ws.Range("I25").Value = 220
$ws.Range("J25").Value = 183
$ws.Range("K25").Value = 20.218579234972
$ws.Range("L25").Value = 51.724137931034
$ws.Range("M25").Value = 6.280193236714

# --- Row 26 (UCR Rape*) ----------------------------------------------------
$ws.Range("F26").Value = "'0"
$ws.Range("G26").Value = "'0"
$ws.Range("H26").Value = "'***.*"
$ws.Range("F23:G23").Copy() | Out-Null
$ws.Range("F26:G26").PasteSpecial(-4122) | Out-Null
$ws.Range("H23").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null

# --- Row 27 (Other Sex Crimes) ---------------------------------------------
$ws.Range("C27").Value = "'0"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = -34.615384615384
$ws.Range("L27").Value = 30.769230769230

# --- Row 28 (Shooting Vic.) --------------------------------------------------
$ws.Range("F28").Value = "'0"
$ws.Range("F23").Copy() | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null
$ws.Range("H28").Value = -100

# --- Row 29 (Shooting Inc.) --------------------------------------------------
$ws.Range("F29").Value = "'0"
$ws.Range("F23").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null
$ws.Range("H29").Value = -100

$excel.CutCopyMode = 0
